# Sync quotation data: row 3 (Dilmatec quote) gets approved, row 5's
# "archived" flag is (re)set to TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quotations")

# --- Row 3: quotation MmE3ZThjNjktMTIyMy00MWQzLTk2N2QtNDc0Nzg4ZDJmYjAxOjU3MDE2 ---
# archived -> TRUE
$ws.Range("F3").Value = $true

# status (display label) -> "Aprovada"
$ws.Range("I3").Value = "Aprovada"

# approvedBy -> filled in now that it has been approved
$ws.Range("L3").Value = "Rafael Machado Barboza"

# approvedSignature -> signature image path for the approval
$ws.Range("N3").Value = "accounts/57016/quotations/2a7e8c69-1223-41d3-967d-474788d2fb01/signatures/51c1ade7-ff37-4ee9-b84d-46ee75381607.png"

# order.id -> cleared (no linked order anymore)
$ws.Range("P3").Value = ""

# status_original (machine status) -> "approved"
$ws.Range("U3").Value = "approved"

# --- Row 5: quotation OTEzYzY1NmQtOWU4OS00Y2VlLWJiNjMtNTI5N2YyMjA1N2M1OjU3MDE2 ---
# archived -> TRUE
$ws.Range("F5").Value = $true
